$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated values for rows 2-7
$ws.Range("A2").Value = "('Acorn Stash', ['Card', '(Place your acorn counters in this area.)'])"
$ws.Range("A3").Value = "('Beeble', ['Token Creature — Beeble', '1/1'])"
$ws.Range("A4").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '4/4'])"
$ws.Range("A5").Value = "('Giant Teddy Bear', ['Token Creature — Giant Teddy Bear', '5/5'])"
$ws.Range("A6").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A7").Value = "('Squirrel', ['Token Creature — Squirrel', '1/1'])"

# Remove now-unused rows 8-20
$ws.Rows("8:20").Delete()
